# Modèle-audit-SEO.xlsx edit script
# Commit message: "toutes images en jpg"
#
# Changes applied (per the OOXML diff):
#  1. Rename the second sheet "Page 2" -> "Page contact"
#  2. Move the selection on sheet 1 to A11 (was A16)
#  3. Highlight (green fill) cells A8, A9, A11 on sheet 1 (category column)
#  4. Replace the old placeholder row 13 (a stray boolean in E13) with a new
#     audit-table row about missing visible link text / Lighthouse recommendation

$wb = $excel.ActiveWorkbook

# --- 1. Rename "Page 2" -> "Page contact" -----------------------------------
$wsContact = $wb.Worksheets.Item(2)
$wsContact.Name = "Page contact"

# --- Work on the main sheet --------------------------------------------------
$ws = $wb.Worksheets.Item(1)

# --- 3. Highlight category cells with the existing green fill ---------------
$greenFill = 5287936   # RGB(0,176,80) == fgColor FF00B050 already used elsewhere
$ws.Range("A8").Interior.Color = $greenFill
$ws.Range("A9").Interior.Color = $greenFill
$ws.Range("A11").Interior.Color = $greenFill

# --- 4. Rebuild row 13 as a new audit entry ----------------------------------
$ws.Range("A13").Value = "Accessibilité"
$ws.Range("B13").Value = "Certains liens n'ont pas de textes visibles"
$ws.Range("C13").Value = "Les technologies d'assistance - comme le lecteur - ne pourront pas naviguer correctement sur ces liens"
$ws.Range("D13").Value = "Ajout d'un texte visible pour chaque lien "
$ws.Range("E13").Value = "OUI"
$ws.Range("F13").Value = "Recommandation Lighthouse"

$ws.Range("A13:D13").Interior.Color = $greenFill
$ws.Range("E13").Interior.Color = $greenFill
$ws.Range("F13").Interior.Color = $greenFill

# --- 2. Move the active selection to A11 -------------------------------------
$ws.Range("A11").Select()
